# Localize the single "Users" sheet into three language-specific sheets:
# English (renamed from "Users"), Հայերեն (Armenian), Русский (Russian).
# Each sheet keeps the same layout (headers in row 1, two data rows below);
# only the header text and the Name/Surname/Gender values are translated.
# Age, Phone Number and Password stay identical across all three sheets.

$wb = $excel.ActiveWorkbook

# --- Step 1: rename the original sheet to "English" -----------------------
$wsEn = $wb.Worksheets.Item(1)
$wsEn.Name = "English"

# --- Step 2: duplicate it (preserves styles/column widths/number formats)--
$wsEn.Copy($null, $wsEn)
$wsHy = $wb.Worksheets.Item(2)
$wsHy.Name = "Հայերեն"

$wsHy.Copy($null, $wsHy)
$wsRu = $wb.Worksheets.Item(3)
$wsRu.Name = "Русский"

# --- Step 3: English sheet - capitalize headers ----------------------------
$wsEn.Cells.Item(1,1).Value = "Name"
$wsEn.Cells.Item(1,2).Value = "Surname"
$wsEn.Cells.Item(1,3).Value = "Gender"
$wsEn.Cells.Item(1,4).Value = "Age"
$wsEn.Cells.Item(1,5).Value = "Phone Number"
$wsEn.Cells.Item(1,6).Value = "Password"
# data rows (name/surname/gender) are unchanged in English - already correct

# --- Step 4: Armenian sheet - translated headers + name/surname/gender ----
$wsHy.Cells.Item(1,1).Value = "Անուն"
$wsHy.Cells.Item(1,2).Value = "Ազգանուն"
$wsHy.Cells.Item(1,3).Value = "Սեռ"
$wsHy.Cells.Item(1,4).Value = "Տարիք"
$wsHy.Cells.Item(1,5).Value = "Հեռ. Համար"
$wsHy.Cells.Item(1,6).Value = "Գաղտնաբառ"

$wsHy.Cells.Item(2,1).Value = "Արամ"
$wsHy.Cells.Item(2,2).Value = "Սուքիասյան"
$wsHy.Cells.Item(2,3).Value = "Արական"

$wsHy.Cells.Item(3,1).Value = "Սուրեն"
$wsHy.Cells.Item(3,2).Value = "Սուրենյան"
$wsHy.Cells.Item(3,3).Value = "Արական"

# --- Step 5: Russian sheet - translated headers + name/surname/gender -----
$wsRu.Cells.Item(1,1).Value = "Имя"
$wsRu.Cells.Item(1,2).Value = "Фамилия"
$wsRu.Cells.Item(1,3).Value = "Пол"
$wsRu.Cells.Item(1,4).Value = "Возраст"
$wsRu.Cells.Item(1,5).Value = "Тел. Номер"
$wsRu.Cells.Item(1,6).Value = "Пароль"

$wsRu.Cells.Item(2,1).Value = "Арам"
$wsRu.Cells.Item(2,2).Value = "Сукиасян"
$wsRu.Cells.Item(2,3).Value = "Мужской"

$wsRu.Cells.Item(3,1).Value = "Сурен"
$wsRu.Cells.Item(3,2).Value = "Суренян"
$wsRu.Cells.Item(3,3).Value = "Мужской"

# --- Step 6: re-select the English sheet as the active tab ----------------
$wsEn.Activate()
